$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance")

# New header column G: another date column
$ws.Range("G1").Value = "July 28 2016"

# Row 2 - update id number, keep name/paid, add new July 28 time
$ws.Range("A2").Value = 123456789
$ws.Range("B2").Value = "xc"
$ws.Range("C2").Value = "PAID"
$ws.Range("E2").Value = "07:53 PM"
$ws.Range("F2").Value = "02:34 PM"
$ws.Range("G2").Value = "03:31 PM"

# Row 3 - update id number, keep the 07:53 PM time in column E
$ws.Range("A3").Value = 246812468
$ws.Range("E3").Value = "07:53 PM"

# Row 4 - update id number, keep name/paid
$ws.Range("A4").Value = 789564123
$ws.Range("B4").Value = "js"
$ws.Range("C4").Value = "PAID"

# Row 5 - update id number, keep name, add new July 28 time
$ws.Range("A5").Value = 258147369
$ws.Range("B5").Value = "ccff"
$ws.Range("G5").Value = "03:31 PM"

# Match the new selection noted in the workbook view
$ws.Range("H10").Select() | Out-Null
